$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.65%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.51%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.086"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.84%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05592"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.15%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.472"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8137"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.04%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8456"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.22%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06977"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.34%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.02817"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.59%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09386"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.19%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.001515"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.21%"
$ws.Range("B13").Value = "One"
$ws.Range("C13").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0005962"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.01%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.006208"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.51%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.608"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.10%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.020"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.01%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.055"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.74%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03204"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.96%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.40%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.773"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.87%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04668"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.42%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001248"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.30%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004565"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "6.47%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009606"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.98%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001938"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03656"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.11%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1367"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.89%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002612"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.31%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003411"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-45.17%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008073"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.16%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005388"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.68%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.44%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002421"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "20.09%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
